$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "deuteron" -> "d" for every data row (target column)
$ws.Range("G2:G11").Value2 = "d"

# Bold + center the header row
$header = $ws.Range("A1:K1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108  # xlCenter

# Move the selection, matching the saved cursor position
$ws.Range("G20").Select() | Out-Null
